# Asesores.xlsx — add a new advisor row for store "ALCAJAM-02":
#   GALLARDO OCAS LUIS FERNANDO
# The new record is inserted right after the existing ALCAJAM-02 block
# (row 119), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at 119 (existing rows 119.. shift down to 120..)
$ws.Rows.Item(119).Insert()

# Populate the new row
$ws.Range("A119").Value = "ALCAJAM-02"
$ws.Range("B119").Value = "GALLARDO OCAS LUIS FERNANDO"

# Restore the author's on-screen viewport/selection at save time
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E103").Select()
